$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row after "Contact" (row 10) for "Jurisdiction"
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Copy formatting from the row below (still has the original style) onto the new row
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Update Version value (now row 3, unaffected by the insert above row 11)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (now row 8, unaffected by the insert above row 11)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
